$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a stray row (old row 13) holding only the professor's name
# ("101761 - Arnaldo Márcio Ramalho Prata") in columns B/C with no label in
# column A. Delete that entire row - this shifts every row below it (old
# rows 14-44) up by one, which also re-aligns the row heights that were
# attached to those rows.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of cells still need their text corrected so the
# labels in column A line up with the right content in columns B/C.
$ws.Range("B10").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C10").Value = "101761 - Arnaldo Márcio Ramalho Prata"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2022" looks like a date, so assigning it with .Value would make
# Excel auto-convert the cell to a date serial number. Copy it from the
# "Ativação:" row instead, which keeps it stored as plain text.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C18").Value = "101761 - Arnaldo Márcio Ramalho Prata"

$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("B20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio"
$ws.Range("C20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio"

$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."
